$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.488.39'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.574.64'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.29'
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3684'
$ws.Range("E7").Value = '  +1.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.07'
$ws.Range("E8").Value = '  -3.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3337'
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.149'
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07580'
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.80'
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.978'
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.950'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").Value = '1.575.88'
$ws.Range("E16").Value = '  +0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001128'
$ws.Range("E17").Value = '  +3.07%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.43'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06744'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.405'
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.65'
$ws.Range("E22").Value = '  +4.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.08'
$ws.Range("E23").Value = '  +1.39%  '
$ws.Range("D24").Value = '22.476.94'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.401'
$ws.Range("E25").Value = '  +0.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.660'
$ws.Range("E26").Value = '  +5.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.77'
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.79'
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.002'
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.38'
$ws.Range("E30").Value = '  +1.96%  '
$ws.Range("D31").Value = '1.751.64'
$ws.Range("E31").Value = '  +0.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.101'
$ws.Range("E32").Value = '  +4.83%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.166'
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.001'
$ws.Range("E34").Value = '  +0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.869'
$ws.Range("E35").Value = '  +3.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08400'
$ws.Range("E36").Value = '  +1.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02486'
$ws.Range("E37").Value = '  +3.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2269'
$ws.Range("E38").Value = '  +2.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06465'
$ws.Range("E39").Value = '  +2.07%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.405'
$ws.Range("E40").Value = '  +2.35%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.300'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.53'
$ws.Range("E42").Value = '  +3.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6346'
$ws.Range("E43").Value = '  +5.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.21'
$ws.Range("E44").Value = '  +4.19%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6186'
$ws.Range("E46").Value = '  +8.91%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.788'
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.084'
$ws.Range("E48").Value = '  +4.09%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '126.37'
$ws.Range("E49").Value = '  +1.93%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.215'
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07237'
$ws.Range("E51").Value = '  -0.37%  '
